# Adds the normalized "course_bucket" sheet and clears residual blank
# inline-string cells left over in "courses" and "buckets" from the old
# wide bucket1..bucket4 schema (schema migration per migrate_schema.py).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Clear leftover blank cells on "courses" (sheet with bucket1..bucket4,
#    prereq columns, etc.) and "buckets" (program/bucket definitions).
#    These cells hold no value (empty inline strings) and should be removed
#    entirely rather than kept as empty placeholders.
# ---------------------------------------------------------------------------
$coursesSheet = $wb.Worksheets.Item("courses")
$blankCoursesCells = @(
    "G2", "M2", "R2", "S2", "T2", "G3", "M3", "R3", "S3", "T3",
    "M4", "R4", "S4", "T4", "G5", "M5", "Q5", "R5", "S5", "T5",
    "G6", "M6", "R6", "S6", "T6", "G7", "M7", "R7", "S7", "T7",
    "M8", "Q8", "R8", "S8", "T8", "M9", "Q9", "R9", "S9", "T9",
    "M10", "Q10", "R10", "S10", "T10", "M11", "Q11", "R11", "S11", "T11",
    "G12", "M12", "R12", "S12", "T12", "M13", "R13", "S13", "T13", "M14",
    "Q14", "R14", "S14", "T14", "M15", "Q15", "R15", "S15", "T15", "R16",
    "S16", "T16", "M17", "M18", "R18", "S18", "T18", "M19", "R19", "S19",
    "T19", "G20", "M20", "R20", "S20", "T20", "G21", "M21", "G22", "M22",
    "R22", "S22", "T22", "G23", "M23", "G24", "M24", "G25", "M25", "G26",
    "M26", "G27", "M27", "G28", "M28", "G29", "M29", "G30", "M30", "G31",
    "M31", "G32", "M32", "M33", "G34", "M34", "R34", "S34", "T34", "G35",
    "M35", "R35", "S35", "T35", "G36", "M36", "R36", "S36", "T36", "G37",
    "M37", "T37", "M39", "G40", "M40", "G41", "M41", "G42", "M42", "T42",
    "M43", "T43", "M44", "T44", "M45", "R45", "S45", "T45"
)
foreach ($addr in $blankCoursesCells) {
    $coursesSheet.Range($addr).ClearContents()
}

$bucketsSheet = $wb.Worksheets.Item("buckets")
$blankBucketsCells = @(
    "F2", "G2", "F3", "F4", "F5", "G5", "E6", "F6", "G6"
)
foreach ($addr in $blankBucketsCells) {
    $bucketsSheet.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Add the new normalized "course_bucket" sheet (track_id, course_code,
#    bucket_id) as the last sheet in the workbook, and populate it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$courseBucketSheet = $wb.Worksheets.Add($null, $lastSheet)
$courseBucketSheet.Name = "course_bucket"

$courseBucketRows = @(
    @("track_id", "course_code", "bucket_id"),
    @("FIN_MAJOR", "ACCO 1030", "BUS_ELEC_4"),
    @("FIN_MAJOR", "ACCO 1031", "BUS_ELEC_4"),
    @("FIN_MAJOR", "BUAD 1001", "BUS_ELEC_4"),
    @("FIN_MAJOR", "ECON 1103", "BUS_ELEC_4"),
    @("FIN_MAJOR", "ECON 1104", "BUS_ELEC_4"),
    @("FIN_MAJOR", "BUAD 1560", "BUS_ELEC_4"),
    @("FIN_MAJOR", "BUAD 2001", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 3001", "CORE"),
    @("FIN_MAJOR", "FINA 3002", "CORE"),
    @("FIN_MAJOR", "FINA 3002", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 3002", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 3002", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 3370", "GENERAL_ELEC"),
    @("FIN_MAJOR", "FINA 3986", "GENERAL_ELEC"),
    @("FIN_MAJOR", "FINA 4001", "CORE"),
    @("FIN_MAJOR", "FINA 4002", "CORE"),
    @("FIN_MAJOR", "FINA 4002", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4002", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4002", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4011", "CORE"),
    @("FIN_MAJOR", "FINA 4020", "CORE"),
    @("FIN_MAJOR", "FINA 4020", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4020", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4020", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4023", "CORE"),
    @("FIN_MAJOR", "FINA 4023", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4023", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4023", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4040", "CORE"),
    @("FIN_MAJOR", "FINA 4040", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4040", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4040", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4050", "CORE"),
    @("FIN_MAJOR", "FINA 4050", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4050", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4050", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4060", "CORE"),
    @("FIN_MAJOR", "FINA 4060", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4060", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4060", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4065", "CORE"),
    @("FIN_MAJOR", "FINA 4065", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4065", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4065", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4075", "CORE"),
    @("FIN_MAJOR", "FINA 4075", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4075", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4075", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4081", "CORE"),
    @("FIN_MAJOR", "FINA 4081", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4081", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4081", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4082", "CORE"),
    @("FIN_MAJOR", "FINA 4082", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4082", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4082", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4084", "CORE"),
    @("FIN_MAJOR", "FINA 4084", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4084", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4084", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4085", "CORE"),
    @("FIN_MAJOR", "FINA 4085", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4085", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4085", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4121", "GENERAL_ELEC"),
    @("FIN_MAJOR", "FINA 4122", "GENERAL_ELEC"),
    @("FIN_MAJOR", "FINA 4123", "GENERAL_ELEC"),
    @("FIN_MAJOR", "FINA 4191", "CORE"),
    @("FIN_MAJOR", "FINA 4191", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4191", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4210", "CORE"),
    @("FIN_MAJOR", "FINA 4210", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4210", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4210", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4211", "CORE"),
    @("FIN_MAJOR", "FINA 4211", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4211", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4211", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4212", "CORE"),
    @("FIN_MAJOR", "FINA 4212", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4212", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4212", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4931", "CORE"),
    @("FIN_MAJOR", "FINA 4931", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4931", "FIN_CHOOSE_1"),
    @("FIN_MAJOR", "FINA 4931", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4953", "CORE"),
    @("FIN_MAJOR", "FINA 4953", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4953", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4986", "CORE"),
    @("FIN_MAJOR", "FINA 4986", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4986", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4989", "CORE"),
    @("FIN_MAJOR", "FINA 4989", "FIN_CHOOSE_2"),
    @("FIN_MAJOR", "FINA 4989", "BUS_ELEC_4"),
    @("FIN_MAJOR", "FINA 4995", "GENERAL_ELEC")
)

$r = 1
foreach ($row in $courseBucketRows) {
    $c = 1
    foreach ($val in $row) {
        $courseBucketSheet.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
